$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
# Column G = "Latest HO Xliff Generate Date"
# Row 3 -> 4ef885aa-2641-4df2-98f1-95cb98f6300f.md
$wsOverview.Range("G3").Value = "2016-08-16 14:15:41"
# Row 4 -> d5bcb239-82c2-4116-a391-477a0cb99245.md (shares the same generated-date string)
$wsOverview.Range("G4").Value = "2016-08-16 14:15:41"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Column E = "Priority" : "ht" -> "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Column H = "Correspond Handoff Datetime"
$wsZhCn.Range("H3").Value = "2016-08-16 14:15:36"
$wsZhCn.Range("H4").Value = "2016-08-16 14:15:36"
# Column K = "Correspond Handback DateTime"
$wsZhCn.Range("K3").Value = "2016-08-16 14:15:53"
$wsZhCn.Range("K4").Value = "2016-08-16 14:15:53"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Column H = "Correspond Handoff Datetime"
$wsDeDe.Range("H3").Value = "2016-08-16 14:15:41"
$wsDeDe.Range("H4").Value = "2016-08-16 14:15:41"
# Column K = "Correspond Handback DateTime"
$wsDeDe.Range("K3").Value = "2016-08-16 14:16:00"
$wsDeDe.Range("K4").Value = "2016-08-16 14:16:00"
